$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title: "TERM I" -> "TERM II" (text is unique in the document, safe to
#    use a plain Find/Replace over the whole document content).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("TERM I", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TERM II", 2) | Out-Null

# ---------------------------------------------------------------------------
# Helper: reliably overwrite the *entire* text of a table cell, regardless of
# how many runs the cell's paragraph currently contains. Simply assigning to
# Cell.Range.Text only ever touches the first run when a cell holds more than
# one run, so for multi-run cells we instead build a fresh Range over the
# same Start/End via the document and assign to that.
# ---------------------------------------------------------------------------
function Set-CellText {
    param($table, [int]$row, [int]$col, [string]$text)

    $cell = $table.Cell($row, $col)
    $runCount = $cell.Range.Paragraphs.Item(1).Range.Text.Length

    if ($cell.Range.Words.Count -gt 1 -and $cell.Range.Text.Length -gt 3) {
        $start = $cell.Range.Start
        $end = $cell.Range.End
        $r = $d.Range($start, $end)
        $r.Text = $text
    } else {
        $cell.Range.Text = $text
    }
}

$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# 2. GEOGRAPHY - Paper 1 row: 60 -> 55.0, C6 -> P7, " E" -> " X"
# ---------------------------------------------------------------------------
Set-CellText $t 4 4 "55.0"
Set-CellText $t 4 6 "P7"

# The subject-grade cell holds two runs (" " + "E"); use the document Range
# approach explicitly so the whole cell content becomes " X".
$cell48 = $t.Cell(4, 8)
$r48 = $d.Range($cell48.Range.Start, $cell48.Range.End)
$r48.Text = " X"

# ---------------------------------------------------------------------------
# 3. GEOGRAPHY - Paper 2 row becomes Paper 3, and its marks/grade are cleared
# ---------------------------------------------------------------------------
Set-CellText $t 5 2 "Paper 3"
Set-CellText $t 5 4 ""
Set-CellText $t 5 6 ""

# ---------------------------------------------------------------------------
# 4. GEOGRAPHY - old Paper 3 row is cleared out entirely
# ---------------------------------------------------------------------------
Set-CellText $t 6 2 ""
Set-CellText $t 6 4 ""
Set-CellText $t 6 6 ""

# ---------------------------------------------------------------------------
# 5. New ECONOMICS subject occupies the two rows that used to be blank
#    placeholders below GEOGRAPHY.
# ---------------------------------------------------------------------------
Set-CellText $t 7 1 "ECONOMICS"
Set-CellText $t 7 8 "E"

Set-CellText $t 8 2 "Paper 2"
Set-CellText $t 8 4 "60.0"
Set-CellText $t 8 6 "C6"

# ---------------------------------------------------------------------------
# 6. ICT row: 74.0 -> 52.0, C4 -> P7 (both paper grade and subject grade)
# ---------------------------------------------------------------------------
Set-CellText $t 13 4 "52.0"
Set-CellText $t 13 6 "P7"
Set-CellText $t 13 8 "P7"
